$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 3.75
$ws.Range("I18").Value = 5
$ws.Range("N18").Value = 1.88
$ws.Range("O18").Value = 1.85
$ws.Range("Q18").Value = 2.72
$ws.Range("R18").Value = 1.87
$ws.Range("S18").Value = 1.85
$ws.Range("U18").Value = 8
$ws.Range("V18").Value = 8.75
$ws.Range("W18").Value = 13
$ws.Range("X18").Value = 14.5
$ws.Range("AA18").Value = 7.6
$ws.Range("AB18").Value = 18
$ws.Range("AC18").Value = 90
$ws.Range("AE18").Value = 12.5
$ws.Range("AG18").Value = 17
$ws.Range("AH18").Value = 100
$ws.Range("AI18").Value = 55
$ws.Range("AJ18").Value = 60
$ws.Range("G19").Value = 2.85
$ws.Range("H19").Value = 3.2
$ws.Range("I19").Value = 2.45
$ws.Range("J19").Value = 1.08
$ws.Range("K19").Value = 7.1
$ws.Range("L19").Value = 1.36
$ws.Range("M19").Value = 3
$ws.Range("N19").Value = 2.05
$ws.Range("O19").Value = 1.72
$ws.Range("P19").Value = 1.42
$ws.Range("Q19").Value = 2.77
$ws.Range("S19").Value = 1.91
$ws.Range("T19").Value = 8
$ws.Range("U19").Value = 15
$ws.Range("W19").Value = 37
$ws.Range("X19").Value = 28
$ws.Range("Z19").Value = 7.1
$ws.Range("AA19").Value = 6.4
$ws.Range("AB19").Value = 15.5
$ws.Range("AE19").Value = 7.7
$ws.Range("AF19").Value = 12.5
$ws.Range("AH19").Value = 28
$ws.Range("G20").Value = 1.88
$ws.Range("I20").Value = 3.85
$ws.Range("U20").Value = 10.5
$ws.Range("X20").Value = 15
$ws.Range("Y20").Value = 24
$ws.Range("AA20").Value = 7.4
$ws.Range("AE20").Value = 12
$ws.Range("AF20").Value = 24
$ws.Range("AI20").Value = 35
$ws.Range("AJ20").Value = 37
$ws.Range("G22").Value = 2.25
$ws.Range("I22").Value = 3.3
$ws.Range("J22").Value = 1.06
$ws.Range("K22").Value = 10
$ws.Range("N22").Value = 2.08
$ws.Range("O22").Value = 1.73
$ws.Range("U22").Value = 10
$ws.Range("W22").Value = 21
$ws.Range("X22").Value = 19
$ws.Range("AD22").Value = 301
$ws.Range("AF22").Value = 17
$ws.Range("AI22").Value = 29
$ws.Range("AJ22").Value = 41
$ws.Range("N24").Value = 2.1
$ws.Range("O24").Value = 1.73
$ws.Range("J25").Value = 1.03
$ws.Range("L25").Value = 1.17
$ws.Range("N25").Value = 1.6
$ws.Range("O25").Value = 2.3
$ws.Range("W25").Value = 13
$ws.Range("G26").Value = 1.98
$ws.Range("H26").Value = 3.05
$ws.Range("I26").Value = 3.65
$ws.Range("N26").Value = 2.07
$ws.Range("O26").Value = 1.6
$ws.Range("Q26").Value = 2.37
$ws.Range("T26").Value = 5.5
$ws.Range("U26").Value = 7.6
$ws.Range("V26").Value = 7.2
$ws.Range("W26").Value = 14.5
$ws.Range("X26").Value = 13.5
$ws.Range("Y26").Value = 24
$ws.Range("AA26").Value = 5.3
$ws.Range("AB26").Value = 13
$ws.Range("AE26").Value = 7.6
$ws.Range("AF26").Value = 15
$ws.Range("AG26").Value = 10.75
$ws.Range("AH26").Value = 40
$ws.Range("AI26").Value = 30
$ws.Range("AJ26").Value = 37
$ws.Range("G27").Value = 5.1
$ws.Range("H27").Value = 3.2
$ws.Range("T27").Value = 9.5
$ws.Range("U27").Value = 23
$ws.Range("V27").Value = 14
$ws.Range("W27").Value = 75
$ws.Range("Z27").Value = 7.3
$ws.Range("AA27").Value = 5.6
$ws.Range("AB27").Value = 15
$ws.Range("AG27").Value = 7.1
$ws.Range("AH27").Value = 10.25
$ws.Range("AJ27").Value = 27
$ws.Range("H29").Value = 4.75
$ws.Range("J29").Value = 1.03
$ws.Range("K29").Value = 17
$ws.Range("R29").Value = 1.73
$ws.Range("S29").Value = 2
$ws.Range("T29").Value = 9.5
$ws.Range("AA29").Value = 9.5
$ws.Range("AB29").Value = 17
$ws.Range("AG29").Value = 23
$ws.Range("G30").Value = 6.5
$ws.Range("H30").Value = 4.75
$ws.Range("I30").Value = 1.42
$ws.Range("R30").Value = 1.62
$ws.Range("S30").Value = 2.2
$ws.Range("T30").Value = 23
$ws.Range("W30").Value = 67
$ws.Range("X30").Value = 41
$ws.Range("Z30").Value = 19
$ws.Range("AA30").Value = 9.5
$ws.Range("AB30").Value = 15
$ws.Range("AE30").Value = 11
$ws.Range("AG30").Value = 9
$ws.Range("AH30").Value = 11
$ws.Range("J31").Value = 1.04
$ws.Range("K31").Value = 13
$ws.Range("N31").Value = 1.8
$ws.Range("O31").Value = 2
$ws.Range("G34").Value = 1.8
$ws.Range("H34").Value = 3.9
$ws.Range("I34").Value = 3.9
$ws.Range("K34").Value = 15
$ws.Range("L34").Value = 1.18
$ws.Range("M34").Value = 4.5
$ws.Range("N34").Value = 1.6
$ws.Range("O34").Value = 2.3
$ws.Range("Z34").Value = 15
$ws.Range("AA34").Value = 7.5
$ws.Range("G35").Value = 7.5
$ws.Range("H35").Value = 5.25
$ws.Range("L35").Value = 1.17
$ws.Range("M35").Value = 5
$ws.Range("N35").Value = 1.57
$ws.Range("O35").Value = 2.35
$ws.Range("T35").Value = 21
$ws.Range("AD35").Value = 251
$ws.Range("AG35").Value = 9
$ws.Range("N40").Value = 1.95
$ws.Range("O40").Value = 1.85
$ws.Range("H42").Value = 11
$ws.Range("I42").Value = 13
$ws.Range("K42").Value = 34
$ws.Range("AB42").Value = 34
$ws.Range("AE53").Value = 10
$ws.Range("AJ53").Value = 35
$ws.Range("L54").Value = 1.4
$ws.Range("T54").Value = 6.5
$ws.Range("U54").Value = 9.5
$ws.Range("AG54").Value = 12
$ws.Range("G55").Value = 1.7
$ws.Range("H55").Value = 3.5
$ws.Range("I55").Value = 5.25
$ws.Range("R55").Value = 2.1
$ws.Range("S55").Value = 1.67
$ws.Range("X55").Value = 15
$ws.Range("AE55").Value = 12
$ws.Range("N57").Value = 1.7
$ws.Range("O57").Value = 2.1
$ws.Range("G59").Value = 2.92
$ws.Range("H59").Value = 2.7
$ws.Range("I59").Value = 2.65
$ws.Range("K59").Value = 4.35
$ws.Range("P59").Value = 1.6
$ws.Range("U59").Value = 13
$ws.Range("V59").Value = 11.5
$ws.Range("W59").Value = 37
$ws.Range("X59").Value = 32
$ws.Range("Y59").Value = 55
$ws.Range("Z59").Value = 4.65
$ws.Range("AA59").Value = 5.6
$ws.Range("AE59").Value = 5.7
$ws.Range("AF59").Value = 11.25
$ws.Range("AG59").Value = 11.25
$ws.Range("AH59").Value = 32
$ws.Range("AI59").Value = 32
$ws.Range("AJ59").Value = 60
